$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1138805618700174
$ws.Range("C2").Value = 0.6694003603824907
$ws.Range("D2").Value = 0.9608326689999721
$ws.Range("E2").Value = 0.9802207246329635
$ws.Range("F2").Value = 0.9843411647153698
$ws.Range("G2").Value = 46
$ws.Range("B3").Value = 0.1339046951553961
$ws.Range("C3").Value = 0.6361935837409688
$ws.Range("D3").Value = 0.871221168554793
$ws.Range("E3").Value = 0.9333922908160283
$ws.Range("F3").Value = 0.9341753966206093
$ws.Range("G3").Value = 45
$ws.Range("B4").Value = 0.1063780764838869
$ws.Range("C4").Value = 0.6141093436083507
$ws.Range("D4").Value = 0.8390184551844652
$ws.Range("E4").Value = 0.9159795058757948
$ws.Range("F4").Value = 0.9202994335924404
$ws.Range("G4").Value = 44
$ws.Range("B5").Value = 0.11775944203778
$ws.Range("C5").Value = 0.5965073955900791
$ws.Range("D5").Value = 0.8555400390622062
$ws.Range("E5").Value = 0.9249540740286548
$ws.Range("F5").Value = 0.9282847517451299
$ws.Range("G5").Value = 43
$ws.Range("B6").Value = 0.1092673841238316
$ws.Range("C6").Value = 0.6065317009167946
$ws.Range("D6").Value = 0.8588038977416446
$ws.Range("E6").Value = 0.9267167300430291
$ws.Range("F6").Value = 0.9314074130616231
$ws.Range("G6").Value = 42
$ws.Range("B7").Value = 0.1420621155324862
$ws.Range("C7").Value = 0.618246900775709
$ws.Range("D7").Value = 0.8875902017685597
$ws.Range("E7").Value = 0.9421200569824207
$ws.Range("F7").Value = 0.9429176904833576
$ws.Range("G7").Value = 41
$ws.Range("B8").Value = 0.1171206258823649
$ws.Range("C8").Value = 0.6341557086697007
$ws.Range("D8").Value = 0.8971434966731959
$ws.Range("E8").Value = 0.947176592126936
$ws.Range("F8").Value = 0.951881405921772
$ws.Range("G8").Value = 40
$ws.Range("B9").Value = 0.1397568759146854
$ws.Range("C9").Value = 0.6260487419841291
$ws.Range("D9").Value = 0.9269108917163016
$ws.Range("E9").Value = 0.9627621158501728
$ws.Range("F9").Value = 0.9650167354246122
$ws.Range("G9").Value = 39
$ws.Range("B10").Value = 0.1262535463842745
$ws.Range("C10").Value = 0.6468408961867226
$ws.Range("D10").Value = 0.9346570565552187
$ws.Range("E10").Value = 0.9667766321934032
$ws.Range("F10").Value = 0.9713636242078157
$ws.Range("G10").Value = 38
$ws.Range("B11").Value = 0.1560923699556315
$ws.Range("C11").Value = 0.6300966329498551
$ws.Range("D11").Value = 0.9461960779220975
$ws.Range("E11").Value = 0.9727261063228937
$ws.Range("F11").Value = 0.9733641012354192
$ws.Range("G11").Value = 37